# Auto-generated update for commit 'Add data for 2022-11-27'
# Updates the 2022 (column I) totals across Citywide Totals, By Neighborhood,
# and each individual neighborhood sheet to reflect one additional day of data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 6637  # I2: 6616 -> 6637
$ws.Cells.Item(3, 9).Value = 6944  # I3: 6919 -> 6944
$ws.Cells.Item(4, 9).Value = 1589  # I4: 1584 -> 1589
$ws.Cells.Item(6, 9).Value = 7992  # I6: 7958 -> 7992
$ws.Cells.Item(7, 9).Value = 23807  # I7: 23722 -> 23807

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(2, 9).Value = 69  # I2: 68 -> 69
$ws.Cells.Item(7, 9).Value = 283  # I7: 282 -> 283

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 9).Value = 63  # I3: 64 -> 63
$ws.Cells.Item(7, 9).Value = 269  # I7: 270 -> 269

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(2, 9).Value = 22  # I2: 21 -> 22
$ws.Cells.Item(6, 9).Value = 22  # I6: 21 -> 22
$ws.Cells.Item(7, 9).Value = 81  # I7: 79 -> 81

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 9).Value = 223  # I2: 221 -> 223
$ws.Cells.Item(3, 9).Value = 245  # I3: 244 -> 245
$ws.Cells.Item(7, 9).Value = 748  # I7: 745 -> 748

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 9).Value = 155  # I3: 153 -> 155
$ws.Cells.Item(6, 9).Value = 108  # I6: 107 -> 108
$ws.Cells.Item(7, 9).Value = 420  # I7: 417 -> 420

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 9).Value = 340  # I3: 338 -> 340
$ws.Cells.Item(6, 9).Value = 276  # I6: 277 -> 276
$ws.Cells.Item(7, 9).Value = 913  # I7: 912 -> 913

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 9).Value = 52  # I6: 51 -> 52
$ws.Cells.Item(7, 9).Value = 211  # I7: 210 -> 211

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(6, 9).Value = 164  # I6: 162 -> 164
$ws.Cells.Item(7, 9).Value = 549  # I7: 547 -> 549

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 9).Value = 97  # I4: 96 -> 97
$ws.Cells.Item(5, 9).Value = 72  # I5: 71 -> 72
$ws.Cells.Item(7, 9).Value = 749  # I7: 748 -> 749
$ws.Cells.Item(8, 9).Value = 1430  # I8: 1419 -> 1430
$ws.Cells.Item(9, 9).Value = 122  # I9: 120 -> 122
$ws.Cells.Item(19, 9).Value = 671  # I19: 670 -> 671
$ws.Cells.Item(20, 9).Value = 590  # I20: 588 -> 590
$ws.Cells.Item(21, 9).Value = 105  # I21: 103 -> 105
$ws.Cells.Item(23, 9).Value = 234  # I23: 233 -> 234
$ws.Cells.Item(24, 9).Value = 65  # I24: 64 -> 65
$ws.Cells.Item(29, 9).Value = 1437  # I29: 1432 -> 1437
$ws.Cells.Item(30, 9).Value = 81  # I30: 79 -> 81
$ws.Cells.Item(33, 9).Value = 1065  # I33: 1062 -> 1065
$ws.Cells.Item(34, 9).Value = 107  # I34: 106 -> 107
$ws.Cells.Item(36, 9).Value = 325  # I36: 324 -> 325
$ws.Cells.Item(37, 9).Value = 748  # I37: 745 -> 748
$ws.Cells.Item(42, 9).Value = 873  # I42: 868 -> 873
$ws.Cells.Item(43, 9).Value = 208  # I43: 206 -> 208
$ws.Cells.Item(46, 9).Value = 53  # I46: 52 -> 53
$ws.Cells.Item(47, 9).Value = 172  # I47: 171 -> 172
$ws.Cells.Item(48, 9).Value = 305  # I48: 303 -> 305
$ws.Cells.Item(49, 9).Value = 158  # I49: 157 -> 158
$ws.Cells.Item(51, 9).Value = 284  # I51: 283 -> 284
$ws.Cells.Item(52, 9).Value = 529  # I52: 525 -> 529
$ws.Cells.Item(53, 9).Value = 260  # I53: 259 -> 260
$ws.Cells.Item(54, 9).Value = 480  # I54: 479 -> 480
$ws.Cells.Item(55, 9).Value = 274  # I55: 273 -> 274
$ws.Cells.Item(63, 9).Value = 71  # I63: 75 -> 71
$ws.Cells.Item(64, 9).Value = 190  # I64: 189 -> 190
$ws.Cells.Item(65, 9).Value = 549  # I65: 547 -> 549
$ws.Cells.Item(67, 9).Value = 913  # I67: 912 -> 913
$ws.Cells.Item(70, 9).Value = 39  # I70: 38 -> 39
$ws.Cells.Item(71, 9).Value = 70  # I71: 68 -> 70
$ws.Cells.Item(72, 9).Value = 95  # I72: 94 -> 95
$ws.Cells.Item(75, 9).Value = 75  # I75: 74 -> 75
$ws.Cells.Item(76, 9).Value = 344  # I76: 343 -> 344
$ws.Cells.Item(78, 9).Value = 320  # I78: 319 -> 320
$ws.Cells.Item(79, 9).Value = 677  # I79: 674 -> 677
$ws.Cells.Item(83, 9).Value = 515  # I83: 512 -> 515
$ws.Cells.Item(84, 9).Value = 211  # I84: 210 -> 211
$ws.Cells.Item(85, 9).Value = 1067  # I85: 1063 -> 1067
$ws.Cells.Item(88, 9).Value = 221  # I88: 220 -> 221
$ws.Cells.Item(89, 9).Value = 283  # I89: 282 -> 283
$ws.Cells.Item(90, 9).Value = 310  # I90: 309 -> 310
$ws.Cells.Item(91, 9).Value = 252  # I91: 251 -> 252
$ws.Cells.Item(94, 9).Value = 243  # I94: 242 -> 243
$ws.Cells.Item(95, 9).Value = 361  # I95: 360 -> 361
$ws.Cells.Item(96, 9).Value = 269  # I96: 270 -> 269
$ws.Cells.Item(97, 9).Value = 194  # I97: 193 -> 194
$ws.Cells.Item(99, 9).Value = 420  # I99: 417 -> 420
$ws.Cells.Item(101, 9).Value = 23807  # I101: 23722 -> 23807

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 174  # I2: 172 -> 174
$ws.Cells.Item(6, 9).Value = 114  # I6: 113 -> 114
$ws.Cells.Item(7, 9).Value = 515  # I7: 512 -> 515

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(6, 9).Value = 75  # I6: 74 -> 75
$ws.Cells.Item(7, 9).Value = 361  # I7: 360 -> 361

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 9).Value = 391  # I3: 390 -> 391
$ws.Cells.Item(6, 9).Value = 344  # I6: 342 -> 344
$ws.Cells.Item(7, 9).Value = 1065  # I7: 1062 -> 1065

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Cells.Item(6, 9).Value = 95  # I6: 94 -> 95
$ws.Cells.Item(7, 9).Value = 158  # I7: 157 -> 158

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 9).Value = 233  # I6: 232 -> 233
$ws.Cells.Item(7, 9).Value = 480  # I7: 479 -> 480

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 423  # I2: 419 -> 423
$ws.Cells.Item(3, 9).Value = 498  # I3: 497 -> 498
$ws.Cells.Item(7, 9).Value = 1437  # I7: 1432 -> 1437

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(4, 9).Value = 24  # I4: 23 -> 24
$ws.Cells.Item(7, 9).Value = 671  # I7: 670 -> 671

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 9).Value = 48  # I2: 47 -> 48
$ws.Cells.Item(6, 9).Value = 156  # I6: 155 -> 156
$ws.Cells.Item(7, 9).Value = 305  # I7: 303 -> 305

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(6, 9).Value = 157  # I6: 156 -> 157
$ws.Cells.Item(7, 9).Value = 344  # I7: 343 -> 344

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 9).Value = 406  # I3: 404 -> 406
$ws.Cells.Item(6, 9).Value = 278  # I6: 276 -> 278
$ws.Cells.Item(7, 9).Value = 1067  # I7: 1063 -> 1067

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 9).Value = 206  # I2: 205 -> 206
$ws.Cells.Item(3, 9).Value = 265  # I3: 263 -> 265
$ws.Cells.Item(4, 9).Value = 57  # I4: 56 -> 57
$ws.Cells.Item(6, 9).Value = 318  # I6: 317 -> 318
$ws.Cells.Item(7, 9).Value = 873  # I7: 868 -> 873

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(3, 9).Value = 80  # I3: 79 -> 80
$ws.Cells.Item(7, 9).Value = 320  # I7: 319 -> 320

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(6, 9).Value = 83  # I6: 82 -> 83
$ws.Cells.Item(7, 9).Value = 274  # I7: 273 -> 274

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(3, 9).Value = 26  # I3: 25 -> 26
$ws.Cells.Item(7, 9).Value = 65  # I7: 64 -> 65

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(6, 9).Value = 18  # I6: 17 -> 18
$ws.Cells.Item(7, 9).Value = 53  # I7: 52 -> 53

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(6, 9).Value = 68  # I6: 67 -> 68
$ws.Cells.Item(7, 9).Value = 234  # I7: 233 -> 234

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(4, 9).Value = 14  # I4: 13 -> 14
$ws.Cells.Item(7, 9).Value = 252  # I7: 251 -> 252

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(2, 9).Value = 8  # I2: 7 -> 8
$ws.Cells.Item(3, 9).Value = 16  # I3: 15 -> 16
$ws.Cells.Item(7, 9).Value = 105  # I7: 103 -> 105

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 9).Value = 196  # I2: 195 -> 196
$ws.Cells.Item(3, 9).Value = 220  # I3: 219 -> 220
$ws.Cells.Item(4, 9).Value = 39  # I4: 38 -> 39
$ws.Cells.Item(7, 9).Value = 677  # I7: 674 -> 677

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 9).Value = 54  # I3: 53 -> 54
$ws.Cells.Item(7, 9).Value = 190  # I7: 189 -> 190

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(4, 9).Value = 38  # I4: 37 -> 38
$ws.Cells.Item(6, 9).Value = 203  # I6: 202 -> 203
$ws.Cells.Item(7, 9).Value = 590  # I7: 588 -> 590

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(4, 9).Value = 12  # I4: 11 -> 12
$ws.Cells.Item(7, 9).Value = 325  # I7: 324 -> 325

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(6, 9).Value = 160  # I6: 156 -> 160
$ws.Cells.Item(7, 9).Value = 529  # I7: 525 -> 529

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(3, 9).Value = 30  # I3: 29 -> 30
$ws.Cells.Item(7, 9).Value = 107  # I7: 106 -> 107

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(6, 9).Value = 140  # I6: 139 -> 140
$ws.Cells.Item(7, 9).Value = 243  # I7: 242 -> 243

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 9).Value = 50  # I3: 49 -> 50
$ws.Cells.Item(7, 9).Value = 172  # I7: 171 -> 172

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(2, 9).Value = 41  # I2: 40 -> 41
$ws.Cells.Item(6, 9).Value = 37  # I6: 36 -> 37
$ws.Cells.Item(7, 9).Value = 122  # I7: 120 -> 122

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(3, 9).Value = 32  # I3: 31 -> 32
$ws.Cells.Item(7, 9).Value = 194  # I7: 193 -> 194

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Cells.Item(6, 9).Value = 8  # I6: 7 -> 8
$ws.Cells.Item(7, 9).Value = 39  # I7: 38 -> 39

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 9).Value = 64  # I2: 63 -> 64
$ws.Cells.Item(7, 9).Value = 221  # I7: 220 -> 221

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 9).Value = 423  # I2: 422 -> 423
$ws.Cells.Item(3, 9).Value = 414  # I3: 408 -> 414
$ws.Cells.Item(6, 9).Value = 459  # I6: 455 -> 459
$ws.Cells.Item(7, 9).Value = 1430  # I7: 1419 -> 1430

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(2, 9).Value = 21  # I2: 20 -> 21
$ws.Cells.Item(7, 9).Value = 72  # I7: 71 -> 72

$ws = $wb.Worksheets.Item('Pullman')
$ws.Cells.Item(3, 9).Value = 23  # I3: 22 -> 23
$ws.Cells.Item(7, 9).Value = 75  # I7: 74 -> 75

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 9).Value = 100  # I2: 99 -> 100
$ws.Cells.Item(7, 9).Value = 310  # I7: 309 -> 310

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(2, 9).Value = 61  # I2: 60 -> 61
$ws.Cells.Item(7, 9).Value = 284  # I7: 283 -> 284

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(2, 9).Value = 41  # I2: 40 -> 41
$ws.Cells.Item(6, 9).Value = 118  # I6: 117 -> 118
$ws.Cells.Item(7, 9).Value = 208  # I7: 206 -> 208

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 9).Value = 121  # I6: 120 -> 121
$ws.Cells.Item(7, 9).Value = 260  # I7: 259 -> 260

$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(3, 9).Value = 26  # I3: 25 -> 26
$ws.Cells.Item(4, 9).Value = 6  # I4: 5 -> 6
$ws.Cells.Item(7, 9).Value = 70  # I7: 68 -> 70

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(6, 9).Value = 46  # I6: 45 -> 46
$ws.Cells.Item(7, 9).Value = 95  # I7: 94 -> 95

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(6, 9).Value = 201  # I6: 200 -> 201
$ws.Cells.Item(7, 9).Value = 749  # I7: 748 -> 749

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(2, 9).Value = 36  # I2: 35 -> 36
$ws.Cells.Item(7, 9).Value = 97  # I7: 96 -> 97
